$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data values for A2:A50 (replacing the previous sample with the
# updated dataset used for the "plots with standard deviation" commit).
$values = @(
    75.675675675675677,
    74.662162162162161,
    75.506756756756758,
    66.047297297297305,
    70.439189189189193,
    70.101351351351354,
    79.391891891891902,
    82.263513513513516,
    80.574324324324323,
    81.25,
    64.695945945945937,
    68.918918918918919,
    73.141891891891902,
    74.155405405405403,
    73.817567567567565,
    64.020270270270274,
    67.736486486486484,
    70.777027027027032,
    76.520270270270274,
    77.027027027027032,
    76.689189189189193,
    63.851351351351347,
    62.5,
    63.006756756756758,
    74.662162162162161,
    69.932432432432435,
    73.648648648648646,
    70.777027027027032,
    67.229729729729726,
    68.581081081081081,
    73.986486486486484,
    74.493243243243242,
    75,
    65.709459459459467,
    67.229729729729726,
    58.108108108108105,
    69.087837837837839,
    65.202702702702695,
    65.540540540540533,
    72.297297297297305,
    72.804054054054063,
    75.168918918918919,
    73.986486486486484,
    75.168918918918919,
    74.155405405405403,
    67.905405405405403,
    65.878378378378372,
    74.324324324324323,
    79.729729729729726
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
